$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) '63.410.59'
Set-TextValue $ws.Cells.Item(2, 5) '  -2.67%  '

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) '3.092.49'
Set-TextValue $ws.Cells.Item(3, 5) '  -1.45%  '

# Row 4
Set-TextValue $ws.Cells.Item(4, 5) '  -0.05%  '

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) '546.96'
Set-TextValue $ws.Cells.Item(5, 5) '  -2.94%  '

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) '138.04'
Set-TextValue $ws.Cells.Item(6, 5) '  -6.33%  '

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) '1.00'
Set-TextValue $ws.Cells.Item(7, 5) '  +0.07%  '

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) '3.080.94'
Set-TextValue $ws.Cells.Item(8, 5) '  -1.50%  '

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) '0.493'
Set-TextValue $ws.Cells.Item(9, 5) '  -0.82%  '

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) '6.58'
Set-TextValue $ws.Cells.Item(10, 5) '  -4.83%  '

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) '0.158'
Set-TextValue $ws.Cells.Item(11, 5) '  +0.57%  '

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) '0.461'
Set-TextValue $ws.Cells.Item(12, 5) '  +0.23%  '

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) '34.97'
Set-TextValue $ws.Cells.Item(13, 5) '  -2.87%  '

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) '0.0000218'
Set-TextValue $ws.Cells.Item(14, 5) '  -1.05%  '

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) '3.574.76'
Set-TextValue $ws.Cells.Item(15, 5) '  -1.80%  '

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) '63.326.76'
Set-TextValue $ws.Cells.Item(16, 5) '  -2.72%  '

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) '3.074.52'
Set-TextValue $ws.Cells.Item(18, 5) '  -1.97%  '

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) '494.76'
Set-TextValue $ws.Cells.Item(19, 5) '  -4.65%  '

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) '6.67'
Set-TextValue $ws.Cells.Item(20, 5) '  -0.75%  '

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) '13.51'
Set-TextValue $ws.Cells.Item(21, 5) '  -2.14%  '

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) '0.703'
Set-TextValue $ws.Cells.Item(22, 5) '  +0.60%  '

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) '7.21'
Set-TextValue $ws.Cells.Item(23, 5) '  -2.93%  '

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) '77.92'
Set-TextValue $ws.Cells.Item(24, 5) '  -0.85%  '

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) '12.28'
Set-TextValue $ws.Cells.Item(25, 5) '  -3.13%  '

# Row 26
Set-TextValue $ws.Cells.Item(26, 5) '  +0.04%  '

# Row 27
Set-TextValue $ws.Cells.Item(27, 2) 'PancakeSwap'
Set-TextValue $ws.Cells.Item(27, 3) 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Cells.Item(27, 4) '2.73'
Set-TextValue $ws.Cells.Item(27, 5) '  -1.89%  '

# Row 28
Set-TextValue $ws.Cells.Item(28, 2) 'RenderToken'
Set-TextValue $ws.Cells.Item(28, 3) 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Cells.Item(28, 4) '8.36'
Set-TextValue $ws.Cells.Item(28, 5) '  -3.30%  '

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) '0.997'
Set-TextValue $ws.Cells.Item(29, 5) '  -0.22%  '

# Row 30
Set-TextValue $ws.Cells.Item(30, 2) 'EthereumClassic'
Set-TextValue $ws.Cells.Item(30, 3) 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Cells.Item(30, 4) '26.58'
Set-TextValue $ws.Cells.Item(30, 5) '  +1.93%  '

# Row 31
Set-TextValue $ws.Cells.Item(31, 2) 'ImmutableX'
Set-TextValue $ws.Cells.Item(31, 3) 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Cells.Item(31, 4) '1.94'
Set-TextValue $ws.Cells.Item(31, 5) '  -8.87%  '

# Row 32
Set-TextValue $ws.Cells.Item(32, 5) '  +0.49%  '

# Row 33
Set-TextValue $ws.Cells.Item(33, 2) 'OKB'
Set-TextValue $ws.Cells.Item(33, 3) 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Cells.Item(33, 4) '61.36'
Set-TextValue $ws.Cells.Item(33, 5) '  +15.95%  '

# Row 34
Set-TextValue $ws.Cells.Item(34, 2) 'Stacks'
Set-TextValue $ws.Cells.Item(34, 3) 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Cells.Item(34, 4) '2.52'
Set-TextValue $ws.Cells.Item(34, 5) '  -6.00%  '

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) '527.70'
Set-TextValue $ws.Cells.Item(35, 5) '  -6.17%  '

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) '5.97'
Set-TextValue $ws.Cells.Item(36, 5) '  -1.34%  '

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) '5.19'
Set-TextValue $ws.Cells.Item(37, 5) '  -4.78%  '

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) '0.0403'
Set-TextValue $ws.Cells.Item(38, 5) '  -7.01%  '

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) '0.0795'
Set-TextValue $ws.Cells.Item(39, 5) '  -2.82%  '

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) '3.056.28'
Set-TextValue $ws.Cells.Item(40, 5) '  -0.62%  '

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) '0.120'
Set-TextValue $ws.Cells.Item(41, 5) '  -1.19%  '

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) '8.13'
Set-TextValue $ws.Cells.Item(42, 5) '  -0.77%  '

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) '2.67'
Set-TextValue $ws.Cells.Item(43, 5) '  -6.79%  '

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) '0.257'
Set-TextValue $ws.Cells.Item(44, 5) '  +0.25%  '

# Row 45
Set-TextValue $ws.Cells.Item(45, 5) '  +0.10%  '

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) '2.05'
Set-TextValue $ws.Cells.Item(46, 5) '  -6.75%  '

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) '121.87'
Set-TextValue $ws.Cells.Item(47, 5) '  +3.22%  '

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) '24.47'
Set-TextValue $ws.Cells.Item(48, 5) '  -1.93%  '

# Row 49
Set-TextValue $ws.Cells.Item(49, 5) '  -0.27%  '

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) '0.0₃0511'
Set-TextValue $ws.Cells.Item(50, 5) '  -2.80%  '

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) '2.35'
Set-TextValue $ws.Cells.Item(51, 5) '  +55.64%  '
